# Scheduled sheet data refresh: updates computed market-price / profit
# columns (currentAveragePrice*, LevePrice*, LeveProfit*) for the rows
# whose underlying item prices changed since the last run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 29: Dripping with Venom
$ws.Range("H29").Value = 4522.1113
$ws.Range("J29").Value = 16950
$ws.Range("L29").Value = 50850
$ws.Range("N29").Value = -51412
# Row 38: Just Give Him a Serum
$ws.Range("H38").Value = 1609.6
$ws.Range("I38").Value = 189.33333
$ws.Range("J38").Value = 3740
$ws.Range("K38").Value = 567.99999
$ws.Range("L38").Value = 11220
$ws.Range("M38").Value = -195.99999
$ws.Range("N38").Value = -11964
# Row 58: A Matter of Vital Importance
$ws.Range("H58").Value = 3466.7144
$ws.Range("J58").Value = 4053.4
$ws.Range("L58").Value = 12160.2
$ws.Range("N58").Value = -12460.2
# Row 62: The Mustache Suits Him
$ws.Range("H62").Value = 1944.8889
$ws.Range("I62").Value = 1684
$ws.Range("J62").Value = 2466.6667
$ws.Range("K62").Value = 1684
$ws.Range("L62").Value = 2466.6667
$ws.Range("M62").Value = -1060
$ws.Range("N62").Value = -3714.6667
# Row 65: Forgery of Convenience (L)
$ws.Range("H65").Value = 1944.8889
$ws.Range("I65").Value = 1684
$ws.Range("J65").Value = 2466.6667
$ws.Range("K65").Value = 8420
$ws.Range("L65").Value = 12333.3335
$ws.Range("M65").Value = -5300
$ws.Range("N65").Value = -18573.3335
# Row 92: Whinier than the Sword
$ws.Range("H92").Value = 13411002
$ws.Range("I92").Value = 2924721.8
$ws.Range("K92").Value = 2924721.8
$ws.Range("M92").Value = -2923473.8
# Row 138: All-night Crafting
$ws.Range("H138").Value = 2858.2754
$ws.Range("I138").Value = 1242.1143
$ws.Range("J138").Value = 4521.9707
$ws.Range("K138").Value = 3726.3429
$ws.Range("L138").Value = 13565.9121
$ws.Range("M138").Value = 1413.6571
$ws.Range("N138").Value = -23845.9121
# Row 139: Something Salty and Ceremonial
$ws.Range("H139").Value = 35000
$ws.Range("J139").Value = 35000
$ws.Range("L139").Value = 35000
$ws.Range("N139").Value = -45280

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust
$ws.Range("H32").Value = 4846
$ws.Range("I32").Value = 4667.263
$ws.Range("J32").Value = 5309.091
$ws.Range("K32").Value = 4667.263
$ws.Range("L32").Value = 5309.091
$ws.Range("M32").Value = -4380.263
$ws.Range("N32").Value = -5883.091
# Row 45: Hollow Hallmarks
$ws.Range("H45").Value = 8324.333000000001
$ws.Range("I45").Value = 10331.6875
$ws.Range("K45").Value = 10331.6875
$ws.Range("M45").Value = -9954.6875
# Row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value = 2210.4312
$ws.Range("I132").Value = 1203.45
$ws.Range("J132").Value = 4448.1665
$ws.Range("K132").Value = 3610.35
$ws.Range("L132").Value = 13344.4995
$ws.Range("M132").Value = -1080.35
$ws.Range("N132").Value = -18404.4995
# Row 139: Backing up My Words
$ws.Range("H139").Value = 40141.125
$ws.Range("J139").Value = 40141.125
$ws.Range("L139").Value = 40141.125
$ws.Range("N139").Value = -50421.125
# Row 141: Essays on Equipment
$ws.Range("H141").Value = 19866.666
$ws.Range("J141").Value = 24800
$ws.Range("L141").Value = 24800
$ws.Range("N141").Value = -35160

$ws = $wb.Worksheets.Item("BSM")
# Row 68: It's All about Execution
$ws.Range("H68").Value = 39000
$ws.Range("J68").Value = 39000
$ws.Range("L68").Value = 39000
$ws.Range("N68").Value = -40622
# Row 71: Too Big to Miss (L)
$ws.Range("H71").Value = 39000
$ws.Range("J71").Value = 39000
$ws.Range("L71").Value = 117000
$ws.Range("N71").Value = -125112
# Row 99: Meddle in Metal
$ws.Range("H99").Value = 250002130
$ws.Range("I99").Value = 1000000000
$ws.Range("J99").Value = 2837
$ws.Range("K99").Value = 1000000000
$ws.Range("L99").Value = 2837
$ws.Range("M99").Value = -999998502
$ws.Range("N99").Value = -5833

$ws = $wb.Worksheets.Item("CRP")
# Row 58: You Do the Heavy Lifting
$ws.Range("H58").Value = 1653.3572
$ws.Range("I58").Value = 975.2941
$ws.Range("K58").Value = 975.2941
$ws.Range("M58").Value = -772.2941
# Row 132: Hull Lotta Damage
$ws.Range("H132").Value = 3625.1304
$ws.Range("I132").Value = 2140.9092
$ws.Range("J132").Value = 4985.6665
$ws.Range("K132").Value = 6422.7276
$ws.Range("L132").Value = 14956.9995
$ws.Range("M132").Value = -3892.7276
$ws.Range("N132").Value = -20016.9995
# Row 136: Turali Quality
$ws.Range("H136").Value = 1653.3572
$ws.Range("I136").Value = 975.2941
$ws.Range("K136").Value = 2925.8823
$ws.Range("M136").Value = -375.8822999999998
# Row 138: Bow Out
$ws.Range("H138").Value = 29020.8
$ws.Range("J138").Value = 29020.8
$ws.Range("L138").Value = 29020.8
$ws.Range("N138").Value = -39300.8

$ws = $wb.Worksheets.Item("CUL")
# Row 49: Leek Soup for the Soul
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("M49").ClearContents()
# Row 68: Such a Butter Face
$ws.Range("H68").Value = 2717.481
$ws.Range("I68").Value = 3963.697
$ws.Range("J68").Value = 1823.4565
$ws.Range("K68").Value = 11891.091
$ws.Range("L68").Value = 5470.3695
$ws.Range("M68").Value = -11080.091
$ws.Range("N68").Value = -7092.3695
# Row 71: No Margarine of Error (L)
$ws.Range("H71").Value = 2717.481
$ws.Range("I71").Value = 3963.697
$ws.Range("J71").Value = 1823.4565
$ws.Range("K71").Value = 35673.273
$ws.Range("L71").Value = 16411.1085
$ws.Range("M71").Value = -31617.273
$ws.Range("N71").Value = -24523.1085
# Row 141: Ocean Explosion
$ws.Range("H141").Value = 14761
$ws.Range("I141").Value = 10610
$ws.Range("K141").Value = 31830
$ws.Range("M141").Value = -26650

$ws = $wb.Worksheets.Item("LTW")
# Row 22: Skin off Their Backs
$ws.Range("H22").Value = 5053023.5
$ws.Range("I22").Value = 27778904
$ws.Range("J22").Value = 2827.7778
$ws.Range("K22").Value = 27778904
$ws.Range("L22").Value = 2827.7778
$ws.Range("M22").Value = -27778609
$ws.Range("N22").Value = -3417.7778
# Row 27: Fire and Hide
$ws.Range("H27").Value = 5053023.5
$ws.Range("I27").Value = 27778904
$ws.Range("J27").Value = 2827.7778
$ws.Range("K27").Value = 27778904
$ws.Range("L27").Value = 2827.7778
$ws.Range("M27").Value = -27778797
$ws.Range("N27").Value = -3041.7778
# Row 68: You Could Say It's a Moving Target
$ws.Range("H68").Value = 100003250
$ws.Range("I68").Value = 3435.5
$ws.Range("J68").Value = 500002500
$ws.Range("K68").Value = 3435.5
$ws.Range("L68").Value = 500002500
$ws.Range("M68").Value = -2686.5
$ws.Range("N68").Value = -500003998
# Row 71: They Call It Bloody Mary (L)
$ws.Range("H71").Value = 100003250
$ws.Range("I71").Value = 3435.5
$ws.Range("J71").Value = 500002500
$ws.Range("K71").Value = 17177.5
$ws.Range("L71").Value = 2500012500
$ws.Range("M71").Value = -13433.5
$ws.Range("N71").Value = -2500019988
# Row 74: Overall, We Blend In
$ws.Range("H74").Value = 20549.25
$ws.Range("I74").Value = 18598.5
$ws.Range("J74").Value = 22500
$ws.Range("K74").Value = 18598.5
$ws.Range("L74").Value = 22500
$ws.Range("M74").Value = -17600.5
$ws.Range("N74").Value = -24496
# Row 77: Eviction Notice (L)
$ws.Range("H77").Value = 20549.25
$ws.Range("I77").Value = 18598.5
$ws.Range("J77").Value = 22500
$ws.Range("K77").Value = 55795.5
$ws.Range("L77").Value = 67500
$ws.Range("M77").Value = -50803.5
$ws.Range("N77").Value = -77484
# Row 132: Tenets of Tanning
$ws.Range("H132").Value = 11717219
$ws.Range("I132").Value = 18062428
$ws.Range("J132").Value = 2985.6155
$ws.Range("K132").Value = 54187284
$ws.Range("L132").Value = 8956.8465
$ws.Range("M132").Value = -54184754
$ws.Range("N132").Value = -14016.8465

$ws = $wb.Worksheets.Item("WVR")
# Row 132: Comfy Cabins
$ws.Range("H132").Value = 1773.1025
$ws.Range("I132").Value = 1624.6296
$ws.Range("J132").Value = 2107.1667
$ws.Range("K132").Value = 4873.8888
$ws.Range("L132").Value = 6321.500100000001
$ws.Range("M132").Value = -2343.8888
$ws.Range("N132").Value = -11381.5001
